# chore: update Sheets via scheduled runner
# Refreshes cached market-board derived figures (price/profit columns H-N)
# on a handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 483500
$ws.Range("I28").Value = 694629.56
$ws.Range("J28").Value = 918.1429
$ws.Range("K28").Value = 694629.56
$ws.Range("L28").Value = 918.1429
$ws.Range("M28").Value = -694144.56
$ws.Range("N28").Value = -1888.1429

$ws.Range("H98").Value = 700754.5
$ws.Range("I98").Value = 932978.5
$ws.Range("K98").Value = 932978.5
$ws.Range("M98").Value = -931480.5

$ws.Range("H122").Value = 700754.5
$ws.Range("I122").Value = 932978.5
$ws.Range("K122").Value = 2798935.5
$ws.Range("M122").Value = -2796485.5

$ws.Range("H129").Value = 982.1177
$ws.Range("I129").Value = 349.81818
$ws.Range("J129").Value = 1156
$ws.Range("K129").Value = 1049.45454
$ws.Range("L129").Value = 3468
$ws.Range("M129").Value = 3950.54546
$ws.Range("N129").Value = -13468

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 254260.86
$ws.Range("I132").Value = 369200.88
$ws.Range("K132").Value = 1107602.64
$ws.Range("M132").Value = -1105072.64

$ws.Range("H4").Value = 300.5
$ws.Range("I4").Value = 300.5
$ws.Range("K4").Value = 300.5
$ws.Range("M4").Value = -184.5

$ws.Range("H23").Value = 17235.295
$ws.Range("J23").Value = 17235.295
$ws.Range("L23").Value = 17235.295
$ws.Range("N23").Value = -17753.295

$ws.Range("H122").Value = 2003
$ws.Range("I122").Value = 2003
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6009
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3559
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 832.3333
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 832.3333
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 832.3333
$ws.Range("N107").Value = -4672.3333
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 484.81818
$ws.Range("I107").Value = 248.57143
$ws.Range("J107").Value = 898.25
$ws.Range("K107").Value = 248.57143
$ws.Range("L107").Value = 898.25
$ws.Range("M107").Value = 1671.42857
$ws.Range("N107").Value = -4738.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 1675
$ws.Range("I116").Value = 850
$ws.Range("K116").Value = 2550
$ws.Range("M116").Value = 892

$ws.Range("H125").Value = 2739.6667
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 2826.65
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 8479.95
$ws.Range("M125").Value = 1920
$ws.Range("N125").Value = -18319.95

$ws.Range("H131").Value = 1566.6207
$ws.Range("I131").Value = 614.75
$ws.Range("J131").Value = 1637.1296
$ws.Range("K131").Value = 1844.25
$ws.Range("L131").Value = 4911.3888
$ws.Range("M131").Value = 3195.75
$ws.Range("N131").Value = -14991.3888

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H102").Value = 1490.1333
$ws.Range("I102").Value = 1267
$ws.Range("J102").Value = 1638.8889
$ws.Range("K102").Value = 1267
$ws.Range("L102").Value = 1638.8889
$ws.Range("M102").Value = 355
$ws.Range("N102").Value = -4882.8889

$ws.Range("H113").Value = 1708.8334
$ws.Range("I113").Value = 1486.6666
$ws.Range("J113").Value = 1931
$ws.Range("K113").Value = 1486.6666
$ws.Range("L113").Value = 1931
$ws.Range("M113").Value = 683.3334
$ws.Range("N113").Value = -6271

$ws.Range("H122").Value = 618980.25
$ws.Range("I122").Value = 1112571.1
$ws.Range("J122").Value = 1991.625
$ws.Range("K122").Value = 3337713.3
$ws.Range("L122").Value = 5974.875
$ws.Range("M122").Value = -3335263.3
$ws.Range("N122").Value = -10874.875

$ws.Range("H132").Value = 3533.7856
$ws.Range("I132").Value = 3122.842
$ws.Range("J132").Value = 4401.3335
$ws.Range("K132").Value = 9368.526
$ws.Range("L132").Value = 13204.0005
$ws.Range("M132").Value = -6838.526
$ws.Range("N132").Value = -18264.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H5").Value = 4166.6665
$ws.Range("I5").Value = 4250
$ws.Range("K5").Value = 4250
$ws.Range("M5").Value = -4137

$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H21").Value = 5006
$ws.Range("I21").Value = 5006
$ws.Range("K21").Value = 5006
$ws.Range("M21").Value = -4832

$ws.Range("H40").Value = 3217.2307
$ws.Range("I40").Value = 2644.3333
$ws.Range("J40").Value = 3520.5293
$ws.Range("K40").Value = 2644.3333
$ws.Range("L40").Value = 3520.5293
$ws.Range("M40").Value = -2508.3333
$ws.Range("N40").Value = -3792.5293

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 3249.75
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10480

$ws.Range("H21").Value = 7000
$ws.Range("I21").Value = 2000
$ws.Range("J21").Value = 12000
$ws.Range("K21").Value = 2000
$ws.Range("L21").Value = 12000
$ws.Range("M21").Value = -1765
$ws.Range("N21").Value = -12470

$ws.Range("H24").Value = 4059.5386
$ws.Range("I24").Value = 1593.3334
$ws.Range("J24").Value = 4799.4
$ws.Range("K24").Value = 1593.3334
$ws.Range("L24").Value = 4799.4
$ws.Range("M24").Value = -1363.3334
$ws.Range("N24").Value = -5259.4

$ws.Range("H28").Value = 29875
$ws.Range("J28").Value = 23166.666
$ws.Range("L28").Value = 23166.666
$ws.Range("N28").Value = -23862.666

$ws.Range("H35").Value = 7000
$ws.Range("I35").Value = 2000
$ws.Range("J35").Value = 12000
$ws.Range("K35").Value = 2000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = -1710
$ws.Range("N35").Value = -12580

$ws.Range("H107").Value = 728.9091
$ws.Range("I107").Value = 897
$ws.Range("J107").Value = 527.2
$ws.Range("K107").Value = 2691
$ws.Range("L107").Value = 1581.6
$ws.Range("M107").Value = -771
$ws.Range("N107").Value = -5421.6
